$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 993
$ws.Range("B2").Value = 963
$ws.Range("C2").Value = 963
$ws.Range("D2").Value = 963
$ws.Range("E2").Value = 991
$ws.Range("F2").Value = 992
$ws.Range("G2").Value = 992
$ws.Range("H2").Value = 999
